$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: clear the existing data block (rows 2-12, cols A-I) ---
# ClearContents keeps the per-cell style (s="...") intact, it just drops
# the value/type, matching the diff's "blank but still styled" cells for
# rows 8-12.
$ws.Range("A2:I12").ClearContents()

# --- Row 2: Fábio / JR Ipatinga ---
$ws.Range("A2").Value = "Fábio"
$ws.Range("B2").Value = "'0730"
$ws.Range("C2").Value = "JR Ipatinga"
$ws.Range("D2").Value = "Algumas câmeras do cliente estão ruins, ele vêm aguardando essa manutenção já tem um tempo."
$ws.Range("G2").Value = "Em andamento"
$ws.Range("H2").Value = "Maxvel: 31 / Forte: 12"

# --- Row 3: Giovani / Mf Eventos ---
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0134"
$ws.Range("C3").Value = "Mf Eventos"
$ws.Range("D3").Value = "Sem comunicação de câmeras, era via DDNS."
$ws.Range("G3").Value = "Pendente"

# --- Row 4: Giovani / RotoPlast ---
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0867"
$ws.Range("C4").Value = "RotoPlast"
$ws.Range("D4").Value = "Sem comunicação de câmeras, era via DDNS."
$ws.Range("G4").Value = "Pendente"

# --- Row 5: Giovani / Carinha de Anjo ---
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0623"
$ws.Range("C5").Value = "Carinha de Anjo"
$ws.Range("D5").Value = "Sem comunicação de câmeras, era via DDNS."
$ws.Range("G5").Value = "Pendente"

# --- Row 6: Giovani / Telemont ---
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'1010"
$ws.Range("C6").Value = "Telemont"
$ws.Range("D6").Value = "Sem comunicação de câmeras, era via DDNS."
$ws.Range("G6").Value = "Pendente"

# --- Row 7: Giovani / Valdemar Amaral ---
$ws.Range("A7").Value = "Giovani"
$ws.Range("B7").Value = "'0840"
$ws.Range("C7").Value = "Valdemar Amaral"
$ws.Range("D7").Value = "Central do cliente parece estar sem bateria."
$ws.Range("G7").Value = "Pendente"

# Rows 8-12 stay blank (already cleared above).

# Row heights on the old rows were explicit (wrapped multi-line text), the
# new/blank content no longer needs that, so autofit back to default.
$ws.Rows("2:12").AutoFit()

# --- Update the saved view/selection state ---
# (topLeftCell / scroll position isn't exposed by this host's window object,
# but the active-cell selection is, and that's what the diff records besides
# the scroll offset.)
$ws.Range("G7").Select()
